$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 102, shifting existing rows 102:110 down to 103:111
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new weekly record
$ws.Range("A102").Value = 9
$ws.Range("B102").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C102").Value = "Metropolitana"
$ws.Range("D102").Value = 45142
$ws.Range("E102").Value = 13
$ws.Range("F102").Value = 100112035
$ws.Range("G102").Value = "Bruselas (repollito)"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 52
$ws.Range("K102").Value = 18000
$ws.Range("L102").Value = 19000
$ws.Range("M102").Value = 18500
$ws.Range("N102").Value = "$/malla 15 kilos"
$ws.Range("O102").Value = "Provincia de Quillota"
$ws.Range("P102").Value = 1233
$ws.Range("Q102").Value = 15
$ws.Range("R102").Value = "Hortaliza"
